# Add the new "2022-Q3" quarterly fund-holdings sheet right after "总计"
# (i.e. immediately before the existing "2022-Q2" sheet), and record the
# new quarter's summary row at the top of the "总计" overview sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new per-fund detail sheet "2022-Q3" before "2022-Q2".
# ---------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("2022-Q2")
$ws = $wb.Worksheets.Add($refSheet)
$ws.Name = "2022-Q3"

# Re-fetch "2022-Q2" (its position shifted once the new sheet was added)
# and reuse its header/row formatting so the new sheet matches its siblings.
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Range("B1:H1").Copy($ws.Range("B1:H1"))
$q2.Range("A2:A10").Copy($ws.Range("A2:A10"))

# Columns B:G hold text values (fund codes/names/figures stored as text,
# matching the rest of the workbook) - force text formatting before
# writing so numeric-looking strings (fund codes, percentages...) keep
# their literal representation (e.g. leading zeros).
$ws.Range("B2:G10").NumberFormat = "@"

$data = @(
  @(0, "561550", "华泰柏瑞中证500增强策略ETF",   "7.81", "99.23", "1.15", "0.0898", 9),
  @(1, "008115", "天弘中证红利低波动100指数C",   "2.44", "94.56", "1.70", "0.0415", 9),
  @(2, "015453", "中欧中证500指数增强A",         "1.20", "88.73", "2.88", "0.0346", 1),
  @(3, "008114", "天弘中证红利低波动100指数A",   "1.89", "94.56", "1.70", "0.0321", 9),
  @(4, "515100", "景顺长城中证红利低波动100ETF", "1.62", "98.63", "1.78", "0.0288", 9),
  @(5, "015454", "中欧中证500指数增强C",         "0.34", "88.73", "2.88", "0.0098", 1),
  @(6, "005966", "安信中证500指数增强C",         "0.16", "92.50", "1.16", "0.0019", 5),
  @(7, "005965", "安信中证500指数增强A",         "0.10", "92.50", "1.16", "0.0012", 5),
  @(8, "005166", "嘉实润和量化6个月定期开放混合", "0.22", "24.64", "0.52", "0.0011", 9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}

# Drop the temporary text-forcing number format so the cells end up with
# the workbook's normal (unstyled) look, same as their neighbour sheets.
$ws.Range("B2:G10").Style = "Normal"

# ---------------------------------------------------------------------
# 2. Update the "总计" overview sheet: insert a new row for 2022-Q3 at
#    the top of the data (row 2), pushing the older quarters down.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# Re-apply the index-column styling (lost by the blank inserted row) by
# copying it from the row immediately below, then clear the formatting
# that Insert() copied down from the header row for B:D.
$summary.Range("A3").Copy($summary.Range("A2"))
$summary.Range("B2:D2").ClearFormats()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 9
$summary.Range("D2").Value = 0.24

# Renumber the auto-increment index column for the rows pushed down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
